$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4810967.5
$ws.Range("I64").Value = 5955381
$ws.Range("J64").Value = 4430
$ws.Range("K64").Value = 5955381
$ws.Range("L64").Value = 4430
$ws.Range("M64").Value = -5955133
$ws.Range("N64").Value = -4926
$ws.Range("H67").Value = 4810967.5
$ws.Range("I67").Value = 5955381
$ws.Range("J67").Value = 4430
$ws.Range("K67").Value = 5955381
$ws.Range("L67").Value = 4430
$ws.Range("M67").Value = -5954523
$ws.Range("N67").Value = -6146
$ws.Range("H76").Value = 3607.5527
$ws.Range("I76").Value = 3002.7407
$ws.Range("J76").Value = 5092.091
$ws.Range("K76").Value = 3002.7407
$ws.Range("L76").Value = 5092.091
$ws.Range("M76").Value = -2687.7407
$ws.Range("N76").Value = -5722.091
$ws.Range("H79").Value = 3607.5527
$ws.Range("I79").Value = 3002.7407
$ws.Range("J79").Value = 5092.091
$ws.Range("K79").Value = 3002.7407
$ws.Range("L79").Value = 5092.091
$ws.Range("M79").Value = -1910.7407
$ws.Range("N79").Value = -7276.091
$ws.Range("H87").Value = 16162
$ws.Range("J87").Value = 16162
$ws.Range("L87").Value = 16162
$ws.Range("N87").Value = -18658
$ws.Range("H90").Value = 16162
$ws.Range("J90").Value = 16162
$ws.Range("L90").Value = 48486
$ws.Range("N90").Value = -60966
$ws.Range("H98").Value = 975.63635
$ws.Range("I98").Value = 1042.6
$ws.Range("J98").Value = 306
$ws.Range("K98").Value = 1042.6
$ws.Range("L98").Value = 306
$ws.Range("M98").Value = 455.4000000000001
$ws.Range("N98").Value = -3302
$ws.Range("H122").Value = 975.63635
$ws.Range("I122").Value = 1042.6
$ws.Range("J122").Value = 306
$ws.Range("K122").Value = 3127.8
$ws.Range("L122").Value = 918
$ws.Range("M122").Value = -677.7999999999997
$ws.Range("N122").Value = -5818

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3462.5
$ws.Range("I63").Value = 3450
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 3450
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -2764
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 3462.5
$ws.Range("I66").Value = 3450
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 17250
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -13818
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 1648.6
$ws.Range("I74").Value = 1210.5518
$ws.Range("J74").Value = 2803.4546
$ws.Range("K74").Value = 1210.5518
$ws.Range("L74").Value = 2803.4546
$ws.Range("M74").Value = -336.5518
$ws.Range("N74").Value = -4551.4546
$ws.Range("H77").Value = 1648.6
$ws.Range("I77").Value = 1210.5518
$ws.Range("J77").Value = 2803.4546
$ws.Range("K77").Value = 6052.759
$ws.Range("L77").Value = 14017.273
$ws.Range("M77").Value = -1684.759
$ws.Range("N77").Value = -22753.273
$ws.Range("H88").Value = 1997.7778
$ws.Range("I88").Value = 1456
$ws.Range("K88").Value = 1456
$ws.Range("M88").Value = -1050
$ws.Range("H91").Value = 1997.7778
$ws.Range("I91").Value = 1456
$ws.Range("K91").Value = 1456
$ws.Range("M91").Value = -52

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3125.9092
$ws.Range("I105").Value = 3058.25
$ws.Range("K105").Value = 3058.25
$ws.Range("M105").Value = -1311.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1171.2727
$ws.Range("I22").Value = 1342.4445
$ws.Range("K22").Value = 1342.4445
$ws.Range("M22").Value = -992.4445000000001
$ws.Range("H62").Value = 2385.7144
$ws.Range("I62").Value = 2385.7144
$ws.Range("K62").Value = 2385.7144
$ws.Range("M62").Value = -1761.7144
$ws.Range("H65").Value = 2385.7144
$ws.Range("I65").Value = 2385.7144
$ws.Range("K65").Value = 11928.572
$ws.Range("M65").Value = -8808.572
$ws.Range("H86").Value = 29416128
$ws.Range("I86").Value = 38465376
$ws.Range("K86").Value = 38465376
$ws.Range("M86").Value = -38464253
$ws.Range("H89").Value = 29416128
$ws.Range("I89").Value = 38465376
$ws.Range("K89").Value = 192326880
$ws.Range("M89").Value = -192321264
$ws.Range("H140").Value = 54897.332
$ws.Range("J140").Value = 54897.332
$ws.Range("L140").Value = 54897.332
$ws.Range("N140").Value = -65257.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3085827.2
$ws.Range("I68").Value = 13333833
$ws.Range("J68").Value = 11425.3
$ws.Range("K68").Value = 40001499
$ws.Range("L68").Value = 34275.89999999999
$ws.Range("M68").Value = -40000688
$ws.Range("N68").Value = -35897.89999999999
$ws.Range("H71").Value = 3085827.2
$ws.Range("I71").Value = 13333833
$ws.Range("J71").Value = 11425.3
$ws.Range("K71").Value = 120004497
$ws.Range("L71").Value = 102827.7
$ws.Range("M71").Value = -120000441
$ws.Range("N71").Value = -110939.7
$ws.Range("H80").Value = 3996.889
$ws.Range("J80").Value = 3871.25
$ws.Range("L80").Value = 11613.75
$ws.Range("N80").Value = -13485.75
$ws.Range("H83").Value = 3996.889
$ws.Range("J83").Value = 3871.25
$ws.Range("L83").Value = 34841.25
$ws.Range("N83").Value = -44201.25
$ws.Range("H127").Value = 1138
$ws.Range("J127").Value = 1138
$ws.Range("L127").Value = 3414
$ws.Range("N127").Value = -13334
$ws.Range("H129").Value = 1637.5667
$ws.Range("I129").Value = 813.4545000000001
$ws.Range("J129").Value = 2114.6843
$ws.Range("K129").Value = 2440.3635
$ws.Range("L129").Value = 6344.0529
$ws.Range("M129").Value = 2559.6365
$ws.Range("N129").Value = -16344.0529
$ws.Range("H131").Value = 858.98
$ws.Range("J131").Value = 883
$ws.Range("L131").Value = 2649
$ws.Range("N131").Value = -12729
$ws.Range("H132").Value = 1502.3529
$ws.Range("I132").Value = 1039.2
$ws.Range("J132").Value = 1695.3334
$ws.Range("K132").Value = 9352.800000000001
$ws.Range("L132").Value = 15258.0006
$ws.Range("M132").Value = -6822.800000000001
$ws.Range("N132").Value = -20318.0006
$ws.Range("H139").Value = 1284.8438
$ws.Range("I139").Value = 855.75
$ws.Range("K139").Value = 2567.25
$ws.Range("M139").Value = 2572.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5339.5806
$ws.Range("I70").Value = 4774.15
$ws.Range("K70").Value = 4774.15
$ws.Range("M70").Value = -4504.15
$ws.Range("H73").Value = 5339.5806
$ws.Range("I73").Value = 4774.15
$ws.Range("K73").Value = 4774.15
$ws.Range("M73").Value = -3838.15

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 169817.5
$ws.Range("I126").Value = 251225
$ws.Range("J126").Value = 7002.5
$ws.Range("K126").Value = 753675
$ws.Range("L126").Value = 21007.5
$ws.Range("M126").Value = -751205
$ws.Range("N126").Value = -25947.5

Write-Host "Update complete"